# Auto-generated edit script for lab2_result.xlsx refactor
$wb = $excel.ActiveWorkbook

# --- Summary sheet: update totals (356/342/14 -> 350/333/17) ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A2").Value = "'350"
$wsSummary.Range("B2").Value = "'333"
$wsSummary.Range("C2").Value = "'17"

# --- All sheet: mark routing-table checks (rows 332-340) as failed ---
$wsAll = $wb.Worksheets.Item("All")
$allUpdates = @(
  @(332, 'False', 'The routing table of as1r1 have the wrong number of routes: 6, expected: 8'),
  @(333, 'False', 'The routing table of as1r2 have the wrong number of routes: 5, expected: 8'),
  @(334, 'False', 'The routing table of as2r1 have the wrong number of routes: 5, expected: 8'),
  @(335, 'False', 'The routing table of as2r2 have the wrong number of routes: 4, expected: 8'),
  @(336, 'False', 'The routing table of as3r1 have the wrong number of routes: 7, expected: 9'),
  @(337, 'False', 'The routing table of root have the wrong number of routes: 1, expected: 2'),
  @(338, 'False', 'The routing table of net have the wrong number of routes: 1, expected: 2'),
  @(339, 'False', 'The routing table of pc have the wrong number of routes: 1, expected: 2'),
  @(340, 'False', 'The routing table of local have the wrong number of routes: 1, expected: 2')
)
foreach ($u in $allUpdates) {
  $r = $u[0]
  $wsAll.Cells.Item($r, 2).Value = "'" + $u[1]
  $wsAll.Cells.Item($r, 3).Value = "'" + $u[2]
}

# Remove the now-obsolete pc.net ping rows (old rows 352-357)
$wsAll.Range("A352:A357").EntireRow.Delete()

# --- Failed sheet: rebuild with the refreshed set of 17 failing checks ---
$wsFailed = $wb.Worksheets.Item("Failed")
$wsFailed.Range("A1:C15").Clear()
$failedRows = @(
  @('Tests Description', 'Passed', 'Reason'),
  @('as2r2 has bgp peer 20.30.1.2', 'False', 'The peering between as2r2 and 20.30.1.2 is not up.'),
  @('as3r2 has bgp peer 20.30.1.1', 'False', 'The session is configured but is in the Active state'),
  @('Checking the routing table of as1r1', 'False', 'The routing table of as1r1 have the wrong number of routes: 6, expected: 8'),
  @('Checking the routing table of as1r2', 'False', 'The routing table of as1r2 have the wrong number of routes: 5, expected: 8'),
  @('Checking the routing table of as2r1', 'False', 'The routing table of as2r1 have the wrong number of routes: 5, expected: 8'),
  @('Checking the routing table of as2r2', 'False', 'The routing table of as2r2 have the wrong number of routes: 4, expected: 8'),
  @('Checking the routing table of as3r1', 'False', 'The routing table of as3r1 have the wrong number of routes: 7, expected: 9'),
  @('Checking the routing table of root', 'False', 'The routing table of root have the wrong number of routes: 1, expected: 2'),
  @('Checking the routing table of net', 'False', 'The routing table of net have the wrong number of routes: 1, expected: 2'),
  @('Checking the routing table of pc', 'False', 'The routing table of pc have the wrong number of routes: 1, expected: 2'),
  @('Checking the routing table of local', 'False', 'The routing table of local have the wrong number of routes: 1, expected: 2'),
  @('Checking that `3.2.0.2` is the local name server for device `as1r1`', 'False', '`resolv.conf` file not found for device `as1r1`'),
  @('Checking that `3.2.0.2` is the local name server for device `as1r2`', 'False', '`resolv.conf` file not found for device `as1r2`'),
  @('Checking that `3.2.0.2` is the local name server for device `as2r1`', 'False', '`resolv.conf` file not found for device `as2r1`'),
  @('Checking that `3.2.0.2` is the local name server for device `as2r2`', 'False', '`resolv.conf` file not found for device `as2r2`'),
  @('Checking that `3.2.0.2` is the local name server for device `as3r1`', 'False', '`resolv.conf` file not found for device `as3r1`'),
  @('Checking that `3.2.0.2` is the local name server for device `as3r2`', 'False', '`resolv.conf` file not found for device `as3r2`')
)
for ($i = 0; $i -lt $failedRows.Length; $i++) {
  $row = $i + 1
  $data = $failedRows[$i]
  $wsFailed.Cells.Item($row, 1).Value = $data[0]
  $wsFailed.Cells.Item($row, 2).Value = "'" + $data[1]
  $wsFailed.Cells.Item($row, 3).Value = "'" + $data[2]
}

Write-Host "edit applied"
